# Refreshed crypto price/volume snapshot (coinranking.com scrape),
# committed by the "Updated symbol list" GitHub Actions job.
#
# D (Price) and E (Volume(1h)) are stored as literal text in the sheet
# (not numbers/percentages), so every numeric-looking replacement is
# entered with a leading apostrophe -- Excel's standard "force text"
# prefix -- which keeps Range.Value a String instead of silently
# coercing it to a Double/Percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '''246.58'),
    @('E2', '''0.96%'),
    @('D3', '''29.83'),
    @('E3', '''9.82%'),
    @('E4', '''1.15%'),
    @('D5', '''0.05704'),
    @('E5', '''0.70%'),
    @('D6', '''6.602'),
    @('E6', '''2.01%'),
    @('D7', '''0.8583'),
    @('E7', '''4.15%'),
    @('D8', '''0.8743'),
    @('E8', '''3.17%'),
    @('E9', '''2.94%'),
    @('D10', '''0.07100'),
    @('E10', '''1.55%'),
    @('D11', '''0.02861'),
    @('E11', '''-0.74%'),
    @('E12', '''-0.05%'),
    @('D13', '''0.001527'),
    @('E13', '''1.13%'),
    @('D14', '''0.04156'),
    @('E14', '''0.34%'),
    @('B15', 'TigerCash'),
    @('C15', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @('D15', '''0.006099'),
    @('E15', '''-2.01%'),
    @('B16', 'UpBots'),
    @('C16', 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'),
    @('D16', '''0.007491'),
    @('E16', '''5,108.37%'),
    @('B17', 'LEO'),
    @('C17', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('D17', '''3.479'),
    @('E17', '''-0.97%'),
    @('B18', 'GateToken'),
    @('C18', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'),
    @('D18', '''3.044'),
    @('E18', '''1.48%'),
    @('B19', 'BTSEToken'),
    @('C19', 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'),
    @('D19', '''2.275'),
    @('E19', '''-1.53%'),
    @('B20', 'One'),
    @('C20', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'),
    @('D20', '''0.01031'),
    @('E20', '''1,629.38%'),
    @('D22', '''0.03268'),
    @('E22', '''3.75%'),
    @('E23', '''3.59%'),
    @('D24', '''2.904'),
    @('E24', '''-18.49%'),
    @('E25', '''0.44%'),
    @('D26', '''0.005085'),
    @('E26', '''14.19%'),
    @('D27', '''0.001219'),
    @('E27', '''0.04%'),
    @('E28', '''23.47%'),
    @('D40', '''0.03749'),
    @('D41', '''0.005680'),
    @('E41', '''-5.95%'),
    @('D42', '''0.1071'),
    @('E42', '''1.77%'),
    @('D43', '''0.002100'),
    @('E43', '''-8.68%'),
    @('D44', '''0.009307'),
    @('E44', '''-3.92%'),
    @('D45', '''0.00005109'),
    @('E45', '''-3.94%'),
    @('E46', '''0.01%'),
    @('D47', '''0.07100'),
    @('E47', '''-29.69%'),
    @('D48', '''0.002712'),
    @('E48', '''5.44%'),
    @('E49', '''0.01%'),
    @('E50', '''0.01%'),
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
